$d = $word.ActiveDocument

# 1. Soziodemographische Daten paragraph
$d.Content.Find.Execute(
    "Die Patientin ist verheiratet, hat einen Sohn (22 Jahre) und eine Tochter (24 Jahre). Sie arbeitet 40 Stunden als Betreuungsassistentin in Teilzeit. Zusammen mit ihrem Sohn, ihrer Tochter und ihrem Ehemann bewohnt sie ein Einfamilienhaus. Der Ehemann ist berentet. Sie hat zwei Geschwister, eine Ärztin (minus 3 Jahre) und einen Krankenpfleger (minus 4 Jahre), zu denen sie kaum Kontakt hat.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Die Patientin ist verheiratet und hat zwei Kinder, einen Sohn (22 Jahre) und eine Tochter (24 Jahre). Sie lebt mit ihrer Familie in einem Einfamilienhaus. Der Ehemann ist berentet. Die Patientin arbeitet aktuell 40 Stunden wöchentlich als Betreuungsassistentin in Teilzeit. Sie hat zwei Geschwister, eine Schwester (3 Jahre jünger, Ärztin) und einen Bruder (4 Jahre jünger, Krankenpfleger), zu denen sie jedoch kaum Kontakt hat.",
    2) | Out-Null

# 2. Symptomatik und psychischer Befund paragraph
$d.Content.Find.Execute(
    "Die Patientin berichtet seit mehreren Jahren unter Schlafstörungen, innerer Unruhe, gereizter Stimmung, Versagensgefühlen, Schuldgefühlen, vermindertem Selbstwert, emotionaler Instabilität, Energieverlust, Konzentrationsstörungen, Verlust der Lebensfreude und Zukunftsängsten. Sie wirkt unsicher und schüchtern, ist aber freundlich und zugewandt. Im formalen Gedankengang besteht Grübeln. Inhaltlich zeigen sich Befürchtungen, den alltäglichen beruflichen Anforderungen nicht mehr gewachsen zu sein. Sie distanziert sich glaubhaft von Suizidalität.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Die Patientin berichtet seit mehreren Jahren unter Schlafstörungen, innerer Unruhe, gereizter Stimmung, Versagensgefühlen, Schuldgefühlen, vermindertem Selbstwert, emotionaler Instabilität, Energieverlust, Konzentrationsstörungen, Verlust der Lebensfreude und Zukunftsängsten. Im Kontakt wirkt sie unsicher und schüchtern, ist aber freundlich und zugewandt. Es zeigt sich Grübeln im formalen Gedankengang. Inhaltlich bestehen Befürchtungen, den alltäglichen beruflichen Anforderungen nicht mehr gewachsen zu sein. Die Gedächtnis- und Aufmerksamkeitsfunktionen sind unauffällig. Es liegen keine Hinweise auf Wahrnehmungs- oder Ich-Störungen vor. Sie raucht nicht und trinkt keinen Alkohol. Sie distanziert sich glaubhaft von Suizidalität.",
    2) | Out-Null

# 3. Somatischer Befund / Konsiliarbericht paragraph
$d.Content.Find.Execute(
    "Somatischer Befund siehe beiliegender ärztlicher Konsiliarbericht. Der Ehemann leidet unter chronischen Kopfschmerzen, Rheuma und hat mehrere Nahrungsmittelunverträglichkeiten. Die Patientin erwartet, dass sie darauf achtet, dass seine Ernährung den Erkrankungen angepasst wird.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Der somatische Befund wird durch einen beiliegenden ärztlichen Konsiliarbericht detailliert dokumentiert. Der Ehemann ist berentet. Die Patientin berichtet keine aktuelle psychopharmakologische Medikation. Vorbehandlungen sind nicht explizit erwähnt, werden aber im Rahmen der Lebensgeschichte thematisiert.",
    2) | Out-Null

# 4. Behandlungsrelevante Angaben zur Lebensgeschichte paragraph
$d.Content.Find.Execute(
    "Die Patientin wuchs in einem Umfeld mit einer dominanten und fordernden Mutter und einem dominanten und egozentrischen, impulsiven und unberechenbaren Vater auf, der seine Frau und Kinder schlug und regelmäßig betrunken nach Hause kam. Sie musste sich um ihre Eltern und Geschwister kümmern und ihre eigenen Bedürfnisse zurückstellen. Sie hat eine internalisierte Erwartung, ihre Rolle als Tochter, Ehefrau, Mutter und Betreuungsassistentin perfekt zu erledigen. Die Beziehung zu ihren Eltern ist belastet, aber sie fühlt sich verpflichtet, bei Notrufen hinzugehen.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Die Patientin wuchs in einem familiären Umfeld mit zwei Geschwistern auf, zu denen sie heute kaum Kontakt hat. Die Beziehung zu ihren Eltern ist belastet, da sie sich bei Notfällen immer noch verpflichtet fühlt, zu helfen, obwohl sie sich zunehmend abgrenzen möchte. Sie internalisiert weiterhin die Erwartung, ihre Rolle als Tochter, Ehefrau, Mutter und Betreuungsassistentin perfekt zu erfüllen. Dies führt zu einer hohen Belastung und erschwert die Abgrenzung. Die Patientin berichtet von einer hohen Arbeitsbelastung und Schwierigkeiten, Nein zu sagen, insbesondere im beruflichen Kontext. Sie hat jedoch begonnen, ihre Grenzen besser wahrzunehmen und sich für ihre Bedürfnisse einzusetzen.",
    2) | Out-Null

# 5. Diagnose zum Zeitpunkt der Antragstellung paragraph
$d.Content.Find.Execute(
    "Die Diagnose zum Zeitpunkt der Antragstellung lautet: (im Text nicht explizit genannt, aber implizit depressive Symptomatik).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Die Diagnose zum Zeitpunkt der Antragstellung lautet F3.1 (Depressive Episode mittleren Schweregrades).",
    2) | Out-Null

# 6. Behandlungsplan und Prognose paragraph
$d.Content.Find.Execute(
    "Ziel der Therapie ist die Stabilisierung der Stimmung, die Steigerung des Selbstwertes, die Förderung der Abgrenzung und die Entwicklung eines funktionierenden Umfelds. Der Behandlungsplan umfasst wöchentliche Sitzungen mit Fokus auf dysfunktionale Gedanken, soziale Kompetenzen und Krisenprophylaxe. Eine Langzeittherapie mit 36 Stunden wird beantragt. Die Prognose ist günstig, da die Patientin motiviert, zuverlässig und reflektiert ist und einen hohen Leidensdruck hat.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ziel der Therapie ist die Reduktion depressiver Symptome, die Steigerung des Selbstwertgefühls und die Verbesserung der Abgrenzungsfähigkeit. Der Behandlungsplan umfasst eine Langzeittherapie mit 36 Stunden, die wöchentlich stattfinden soll. Im Rahmen der Therapie werden positive Aktivitäten aufgebaut, dysfunktionale Gedankenmuster bearbeitet und soziale Kompetenzen erweitert. Ein individueller Krisenplan soll erstellt werden, um Rückfälle zu vermeiden. Die Prognose wird als günstig eingeschätzt, da die Patientin motiviert ist, regelmäßig zu den Terminen kommt, ihre Hausaufgaben erledigt und einen hohen Leidensdruck verspürt. Sie zeigt Bereitschaft zur Veränderung und wünscht sich ein funktionierendes Umfeld, in dem sie ihre Bedürfnisse äußern und umsetzen kann.",
    2) | Out-Null

# 7. Closing paragraph: append two more line-break + text runs after
#    "[Name des Therapeuten/der Therapeutin]"
$vbreak = [char]11
$closing = $d.Paragraphs.Last
$closingRange = $closing.Range
$closingRange.InsertAfter($vbreak + "[Berufsbezeichnung]" + $vbreak + "[Kontaktdaten]")

Write-Output "done"
